# 422 error with polygon search
# Truncate the sqm (column D) values to whole numbers for rows 2-4,
# and update the active selection to D7 (as if the user had just
# finished editing column D and moved down/right).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 259
$ws.Range("D3").Value = 548
$ws.Range("D4").Value = 134

$ws.Range("D7").Select()
